$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value2 = '''305.76'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value2 = '''0.81%'
$ws.Range('E2').Style = 'Normal'
$ws.Range('G2').Value2 = '''11'
$ws.Range('G2').Style = 'Normal'
$ws.Range('D3').Value2 = '''36.07'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value2 = '''-2.90%'
$ws.Range('E3').Style = 'Normal'
$ws.Range('G3').Value2 = '''11'
$ws.Range('G3').Style = 'Normal'
$ws.Range('E4').Value2 = '''2.16%'
$ws.Range('E4').Style = 'Normal'
$ws.Range('G4').Value2 = '''11'
$ws.Range('G4').Style = 'Normal'
$ws.Range('D5').Value2 = '''0.07874'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value2 = '''0.62%'
$ws.Range('E5').Style = 'Normal'
$ws.Range('G5').Value2 = '''11'
$ws.Range('G5').Style = 'Normal'
$ws.Range('D6').Value2 = '''2.132'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value2 = '''-3.56%'
$ws.Range('E6').Style = 'Normal'
$ws.Range('G6').Value2 = '''11'
$ws.Range('G6').Style = 'Normal'
$ws.Range('D7').Value2 = '''7.930'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value2 = '''-1.16%'
$ws.Range('E7').Style = 'Normal'
$ws.Range('G7').Value2 = '''11'
$ws.Range('G7').Style = 'Normal'
$ws.Range('B8').Value2 = '''GateToken'
$ws.Range('B8').Style = 'Normal'
$ws.Range('C8').Value2 = '''https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('C8').Style = 'Normal'
$ws.Range('D8').Value2 = '''4.114'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value2 = '''1.86%'
$ws.Range('E8').Style = 'Normal'
$ws.Range('G8').Value2 = '''11'
$ws.Range('G8').Style = 'Normal'
$ws.Range('B9').Value2 = '''MXToken'
$ws.Range('B9').Style = 'Normal'
$ws.Range('C9').Value2 = '''https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('C9').Style = 'Normal'
$ws.Range('D9').Value2 = '''0.9244'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value2 = '''1.07%'
$ws.Range('E9').Style = 'Normal'
$ws.Range('G9').Value2 = '''11'
$ws.Range('G9').Style = 'Normal'
$ws.Range('B10').Value2 = '''LiechtensteinCryptoassetsExchange'
$ws.Range('B10').Style = 'Normal'
$ws.Range('C10').Value2 = '''https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('C10').Style = 'Normal'
$ws.Range('D10').Value2 = '''0.09692'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value2 = '''-0.04%'
$ws.Range('E10').Style = 'Normal'
$ws.Range('G10').Value2 = '''11'
$ws.Range('G10').Style = 'Normal'
$ws.Range('B11').Value2 = '''WazirX'
$ws.Range('B11').Style = 'Normal'
$ws.Range('C11').Value2 = '''https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('C11').Style = 'Normal'
$ws.Range('D11').Value2 = '''0.1860'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value2 = '''-1.22%'
$ws.Range('E11').Style = 'Normal'
$ws.Range('G11').Value2 = '''11'
$ws.Range('G11').Style = 'Normal'
$ws.Range('B12').Value2 = '''MandalaExchangeToken'
$ws.Range('B12').Style = 'Normal'
$ws.Range('C12').Value2 = '''https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('C12').Style = 'Normal'
$ws.Range('D12').Value2 = '''0.08701'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value2 = '''1.67%'
$ws.Range('E12').Style = 'Normal'
$ws.Range('G12').Value2 = '''11'
$ws.Range('G12').Style = 'Normal'
$ws.Range('B13').Value2 = '''BitrueCoin'
$ws.Range('B13').Style = 'Normal'
$ws.Range('C13').Value2 = '''https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('C13').Style = 'Normal'
$ws.Range('D13').Value2 = '''0.03557'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value2 = '''-0.45%'
$ws.Range('E13').Style = 'Normal'
$ws.Range('G13').Value2 = '''11'
$ws.Range('G13').Style = 'Normal'
$ws.Range('B14').Value2 = '''BitMartToken'
$ws.Range('B14').Style = 'Normal'
$ws.Range('C14').Value2 = '''https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('C14').Style = 'Normal'
$ws.Range('D14').Value2 = '''0.09944'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value2 = '''-0.19%'
$ws.Range('E14').Style = 'Normal'
$ws.Range('G14').Value2 = '''11'
$ws.Range('G14').Style = 'Normal'
$ws.Range('B15').Value2 = '''BitForexToken'
$ws.Range('B15').Style = 'Normal'
$ws.Range('C15').Value2 = '''https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('C15').Style = 'Normal'
$ws.Range('D15').Value2 = '''0.001437'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value2 = '''-3.05%'
$ws.Range('E15').Style = 'Normal'
$ws.Range('G15').Value2 = '''11'
$ws.Range('G15').Style = 'Normal'
$ws.Range('B16').Value2 = '''TigerCash'
$ws.Range('B16').Style = 'Normal'
$ws.Range('C16').Value2 = '''https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('C16').Style = 'Normal'
$ws.Range('D16').Value2 = '''0.005628'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value2 = '''-0.37%'
$ws.Range('E16').Style = 'Normal'
$ws.Range('G16').Value2 = '''11'
$ws.Range('G16').Style = 'Normal'
$ws.Range('B17').Value2 = '''LEO'
$ws.Range('B17').Style = 'Normal'
$ws.Range('C17').Value2 = '''https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('C17').Style = 'Normal'
$ws.Range('D17').Value2 = '''3.450'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value2 = '''-0.33%'
$ws.Range('E17').Style = 'Normal'
$ws.Range('G17').Value2 = '''11'
$ws.Range('G17').Style = 'Normal'
$ws.Range('D18').Value2 = '''2.774'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value2 = '''22.59%'
$ws.Range('E18').Style = 'Normal'
$ws.Range('G18').Value2 = '''11'
$ws.Range('G18').Style = 'Normal'
$ws.Range('D19').Value2 = '''0.3396'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value2 = '''-1.92%'
$ws.Range('E19').Style = 'Normal'
$ws.Range('G19').Value2 = '''11'
$ws.Range('G19').Style = 'Normal'
$ws.Range('D20').Value2 = '''0.1326'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value2 = '''1.98%'
$ws.Range('E20').Style = 'Normal'
$ws.Range('G20').Value2 = '''11'
$ws.Range('G20').Style = 'Normal'
$ws.Range('D21').Value2 = '''5.182'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value2 = '''7.93%'
$ws.Range('E21').Style = 'Normal'
$ws.Range('G21').Value2 = '''11'
$ws.Range('G21').Style = 'Normal'
$ws.Range('D22').Value2 = '''0.1998'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value2 = '''-12.96%'
$ws.Range('E22').Style = 'Normal'
$ws.Range('G22').Value2 = '''11'
$ws.Range('G22').Style = 'Normal'
$ws.Range('D23').Value2 = '''0.04561'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value2 = '''-1.17%'
$ws.Range('E23').Style = 'Normal'
$ws.Range('G23').Value2 = '''11'
$ws.Range('G23').Style = 'Normal'
$ws.Range('D24').Value2 = '''0.005044'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value2 = '''5.39%'
$ws.Range('E24').Style = 'Normal'
$ws.Range('G24').Value2 = '''11'
$ws.Range('G24').Style = 'Normal'
$ws.Range('D25').Value2 = '''0.001236'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value2 = '''0.32%'
$ws.Range('E25').Style = 'Normal'
$ws.Range('G25').Value2 = '''11'
$ws.Range('G25').Style = 'Normal'
$ws.Range('G26').Value2 = '''11'
$ws.Range('G26').Style = 'Normal'
$ws.Range('D27').Value2 = '''0.0004748'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value2 = '''-0.07%'
$ws.Range('E27').Style = 'Normal'
$ws.Range('G27').Value2 = '''11'
$ws.Range('G27').Style = 'Normal'
$ws.Range('G28').Value2 = '''11'
$ws.Range('G28').Style = 'Normal'
$ws.Range('G29').Value2 = '''11'
$ws.Range('G29').Style = 'Normal'
$ws.Range('G30').Value2 = '''11'
$ws.Range('G30').Style = 'Normal'
$ws.Range('G31').Value2 = '''11'
$ws.Range('G31').Style = 'Normal'
$ws.Range('G32').Value2 = '''11'
$ws.Range('G32').Style = 'Normal'
$ws.Range('G33').Value2 = '''11'
$ws.Range('G33').Style = 'Normal'
$ws.Range('G34').Value2 = '''11'
$ws.Range('G34').Style = 'Normal'
$ws.Range('G35').Value2 = '''11'
$ws.Range('G35').Style = 'Normal'
$ws.Range('G36').Value2 = '''11'
$ws.Range('G36').Style = 'Normal'
$ws.Range('G37').Value2 = '''11'
$ws.Range('G37').Style = 'Normal'
$ws.Range('G38').Value2 = '''11'
$ws.Range('G38').Style = 'Normal'
$ws.Range('D39').Value2 = '''0.01850'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value2 = '''4.03%'
$ws.Range('E39').Style = 'Normal'
$ws.Range('G39').Value2 = '''11'
$ws.Range('G39').Style = 'Normal'
$ws.Range('D40').Value2 = '''0.04781'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value2 = '''0.75%'
$ws.Range('E40').Style = 'Normal'
$ws.Range('G40').Value2 = '''11'
$ws.Range('G40').Style = 'Normal'
$ws.Range('D41').Value2 = '''0.007535'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value2 = '''-6.45%'
$ws.Range('E41').Style = 'Normal'
$ws.Range('G41').Value2 = '''11'
$ws.Range('G41').Style = 'Normal'
$ws.Range('D42').Value2 = '''0.1399'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value2 = '''0.62%'
$ws.Range('E42').Style = 'Normal'
$ws.Range('G42').Value2 = '''11'
$ws.Range('G42').Style = 'Normal'
$ws.Range('D43').Value2 = '''0.007743'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value2 = '''1.23%'
$ws.Range('E43').Style = 'Normal'
$ws.Range('G43').Value2 = '''11'
$ws.Range('G43').Style = 'Normal'
$ws.Range('D44').Value2 = '''0.002229'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value2 = '''3.14%'
$ws.Range('E44').Style = 'Normal'
$ws.Range('G44').Value2 = '''11'
$ws.Range('G44').Style = 'Normal'
$ws.Range('D45').Value2 = '''0.01132'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value2 = '''14.04%'
$ws.Range('E45').Style = 'Normal'
$ws.Range('G45').Value2 = '''11'
$ws.Range('G45').Style = 'Normal'
$ws.Range('D46').Value2 = '''0.00006322'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value2 = '''3.55%'
$ws.Range('E46').Style = 'Normal'
$ws.Range('G46').Value2 = '''11'
$ws.Range('G46').Style = 'Normal'
$ws.Range('D47').Value2 = '''0.00000000750'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value2 = '''-0.11%'
$ws.Range('E47').Style = 'Normal'
$ws.Range('G47').Value2 = '''11'
$ws.Range('G47').Style = 'Normal'
$ws.Range('D48').Value2 = '''0.0005798'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value2 = '''-0.05%'
$ws.Range('E48').Style = 'Normal'
$ws.Range('G48').Value2 = '''11'
$ws.Range('G48').Style = 'Normal'
$ws.Range('D49').Value2 = '''47.59'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value2 = '''521.67%'
$ws.Range('E49').Style = 'Normal'
$ws.Range('G49').Value2 = '''11'
$ws.Range('G49').Style = 'Normal'
$ws.Range('D50').Value2 = '''0.001999'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value2 = '''-25.69%'
$ws.Range('E50').Style = 'Normal'
$ws.Range('G50').Value2 = '''11'
$ws.Range('G50').Style = 'Normal'
$ws.Range('D51').Value2 = '''0.00002099'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value2 = '''-0.11%'
$ws.Range('E51').Style = 'Normal'
$ws.Range('G51').Value2 = '''11'
$ws.Range('G51').Style = 'Normal'
